$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093962788581848
$ws.Range("B1").Value = 1.049266695976257
$ws.Range("C1").Value = 1.026624917984009
$ws.Range("D1").Value = 1.252538204193115
$ws.Range("E1").Value = 1.166861414909363
